$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("startup")

# Replace the Cypher query for the "ParticipantsTab" row (row 2, column B /
# "query") with the corrected/expanded version.
$newQuery = @'
MATCH (p:participant)-->(s:study)
OPTIONAL MATCH (samp:sample)-->(p)
OPTIONAL MATCH (p)<--(diag:diagnosis)
OPTIONAL MATCH (samp)<--(f:file)
OPTIONAL MATCH (f)<--(g:genomic_info)
WITH s, p, samp, f, g, diag
WHERE s.study_name in ["Discovery of Colorectal Cancer Susceptibility Genes in High-Risk Families"]
OPTIONAL MATCH (p)-->(s:study)
OPTIONAL MATCH (samp:sample)-->(p)
WITH s, p, apoc.coll.sort(collect(distinct samp.sample_id)) as samp
RETURN 
coalesce(p.participant_id,'') as `Participant ID`,
coalesce(s.study_name, '') as `Study Name`,
coalesce(s.phs_accession,'') as `Accession`,
coalesce(p.gender,'') as `Gender`,
coalesce(apoc.text.join(samp, ','), '') as `Samples`
ORDER BY p.participant_id
LIMIT 100
'@

$ws.Range("B2").Value = $newQuery

# The longer query text wraps onto more lines; match Excel's recalculated
# row height for row 2.
$ws.Rows.Item(2).RowHeight = 299.25

# Update the active selection to B2, matching the saved view state.
$ws.Range("B2").Select()
